$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "StatQuery" text (column C, rows 2-5) gains a two-space indent before the
# WHERE clause and a trailing space after ['Melanoma'].
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
  WHERE diag.disease_term IN ['Melanoma'] 
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Preserve the existing (autofit) row heights - only the text content is
# changing, not the desired rendered row height.
$h2 = $ws.Rows.Item(2).RowHeight
$h3 = $ws.Rows.Item(3).RowHeight
$h4 = $ws.Rows.Item(4).RowHeight
$h5 = $ws.Rows.Item(5).RowHeight

$ws.Range("C2").Value2 = $statQuery
$ws.Range("C3").Value2 = $statQuery
$ws.Range("C4").Value2 = $statQuery
$ws.Range("C5").Value2 = $statQuery

$ws.Rows.Item(2).RowHeight = $h2
$ws.Rows.Item(3).RowHeight = $h3
$ws.Rows.Item(4).RowHeight = $h4
$ws.Rows.Item(5).RowHeight = $h5

# Move the active selection from B5 to B2.
[void]$ws.Range("B2").Select()
